$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.693.89'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '2.654.36'
$ws.Range("E3").Value = '  -0.53%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '596.46'
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").Value = '156.85'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").Value = '0.655'
$ws.Range("E7").Value = '  +5.05%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '0.125'
$ws.Range("E9").Value = '  -3.66%  '
$ws.Range("D10").Value = "'" + '0.400'
$ws.Range("E10").Value = '  -0.81%  '
$ws.Range("D11").Value = '5.82'
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("E12").Value = '  +1.46%  '
$ws.Range("D13").Value = '28.72'
$ws.Range("E13").Value = '  -1.94%  '
$ws.Range("D14").Value = "'" + '0.0000191'
$ws.Range("E14").Value = '  -3.84%  '
$ws.Range("D15").Value = '3.138.95'
$ws.Range("E15").Value = '  -0.31%  '
$ws.Range("D16").Value = '65.616.47'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '2.689.26'
$ws.Range("E17").Value = '  +0.64%  '
$ws.Range("D18").Value = '12.53'
$ws.Range("E18").Value = '  -2.11%  '
$ws.Range("D19").Value = '4.77'
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("D20").Value = '348.75'
$ws.Range("E20").Value = '  -0.37%  '
$ws.Range("D21").Value = '7.39'
$ws.Range("E21").Value = '  -2.91%  '
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").Value = '69.77'
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").Value = '1.81'
$ws.Range("E24").Value = '  +9.85%  '
$ws.Range("D25").Value = "'" + '0.0000111'
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("D26").Value = '9.49'
$ws.Range("E26").Value = '  -1.78%  '
$ws.Range("D27").Value = '1.62'
$ws.Range("E27").Value = '  +2.54%  '
$ws.Range("D28").Value = '565.31'
$ws.Range("E28").Value = '  +7.05%  '
$ws.Range("D29").Value = '8.03'
$ws.Range("E29").Value = '  -0.69%  '
$ws.Range("E30").Value = '  -2.87%  '
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.13%  '
$ws.Range("D32").Value = '2.12'
$ws.Range("E32").Value = '  -1.25%  '
$ws.Range("D33").Value = '1.81'
$ws.Range("E33").Value = '  +3.22%  '
$ws.Range("D34").Value = '6.65'
$ws.Range("E34").Value = '  +3.10%  '
$ws.Range("D35").Value = '5.43'
$ws.Range("E35").Value = '  -1.28%  '
$ws.Range("D36").Value = "'" + '0.420'
$ws.Range("E36").Value = '  -0.60%  '
$ws.Range("D37").Value = '20.47'
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '1.92'
$ws.Range("E39").Value = '  -0.72%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = '154.71'
$ws.Range("E40").Value = '  -2.38%  '
$ws.Range("D41").Value = '160.24'
$ws.Range("E41").Value = '  -2.52%  '
$ws.Range("D42").Value = '4.07'
$ws.Range("E42").Value = '  -2.12%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").Value = '2.27'
$ws.Range("E44").Value = '  -1.93%  '
$ws.Range("D45").Value = '22.66'
$ws.Range("E45").Value = '  -0.85%  '
$ws.Range("D46").Value = '0.638'
$ws.Range("E46").Value = '  -0.48%  '
$ws.Range("E47").Value = '  +1.36%  '
$ws.Range("D48").Value = '0.0253'
$ws.Range("E48").Value = '  -1.94%  '
$ws.Range("D49").Value = '19.71'
$ws.Range("E49").Value = '  -2.41%  '
$ws.Range("D50").Value = '0.0₆0243'
$ws.Range("E50").Value = '  -0.25%  '
$ws.Range("D51").Value = "'" + '0.800'
$ws.Range("E51").Value = '  -2.60%  '
